# Append a new row to the "ランサーズ" (Lancers) sheet and refresh the
# "取得日時" (fetched-at) timestamp on all existing data rows.
#
# Source diff summary:
#   - A2:A19 timestamp "2025-10-15 01:18:04" -> "2025-10-15 01:45:43"
#   - new row 20 appended with a fresh listing + its hyperlink in column F
#   - dimension / hyperlinks collection grow accordingly (handled automatically)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-15 01:45:43"

# Refresh the "取得日時" column for every existing data row (2-19).
for ($i = 2; $i -le 19; $i++) {
    $ws.Cells.Item($i, 1).Value = $newTimestamp
}

# Append the new listing as row 20.
$row = 20
$ws.Cells.Item($row, 1).Value = $newTimestamp
$ws.Cells.Item($row, 2).Value = "Access 32bitから64bitへの修正改善依頼"
$ws.Cells.Item($row, 3).Value = "システム開発"
$ws.Cells.Item($row, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item($row, 5).Value = "期限情報なし"
$ws.Cells.Item($row, 6).Value = "https://www.lancers.jp/work/detail/5413333"
$ws.Cells.Item($row, 7).Value = 10

# Wire up the hyperlink on F20, matching the style used by the rest of
# column F (the Hyperlink cell style).
$ws.Hyperlinks.Add($ws.Cells.Item($row, 6), "https://www.lancers.jp/work/detail/5413333")
$ws.Cells.Item($row, 6).Style = "Hyperlink"
